# Atualização automática do relatório BI
#
# The report's "Tempo total na fase <X> (dias)" columns (AK/AN/AQ/AT/AW) hold
# the elapsed time a row has spent in its CURRENT phase, measured up to the
# moment the report was generated. Re-running the report later moves that
# "as of" instant forward, so every row that is still sitting in an open
# phase (no exit timestamp recorded yet) gets its elapsed-time figure bumped
# by the same amount of wall-clock time that passed between generations.
#
# This recomputes those live "Tempo total" values for the current sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# How far the report's "as of now" reference point advanced since the last
# generation (in days) - derived from the BI refresh that produced this diff.
$delta = 0.16669

# Phase name (as shown in column C, "Fase atual") -> the "Tempo total" column
# that tracks time-in-phase for that stage of the workflow.
$phaseColumn = @{
    "Backlog"                    = "AK"
    "Construção no Canvas"       = "AN"
    "Validação"                  = "AQ"
    "Publicar na plataforma"     = "AT"
    "Concluído"                  = "AW"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $phase = $ws.Cells.Item($row, 3).Value2
    if ($phase -eq $null) { continue }

    $col = $phaseColumn[$phase]
    if ($col -eq $null) { continue }

    $cell = $ws.Range("$col$row")
    $current = $cell.Value2
    if ($current -eq $null) { continue }

    $cell.Value = $current + $delta
}
